# Fruta / hortaliza, semanal
# Update the weekly price data table: dates and price/volume figures for
# rows 2-15 are refreshed to reflect the latest pull of market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad de comercializacion),
# O (Origen), P (Precio $/Kg) and Q (Kg o Unidades).
$rows = @{
    2  = @{ D = 44200; J = 10; K = 9000;  L = 9000;  M = 9000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 150; Q = 60 }
    3  = @{ D = 44284; J = 35; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 167; Q = 60 }
    4  = @{ D = 44277; J = 25; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 167; Q = 60 }
    5  = @{ D = 44291; J = 20; K = 9000;  L = 9000;  M = 9000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 150; Q = 60 }
    6  = @{ D = 44243; J = 80; K = 10000; L = 11000; M = 10375; N = '$/caja 60 unidades'; O = 'Provincia de Quillota';        P = 173; Q = 60 }
    8  = @{ D = 44585; J = 30; K = 11000; L = 11000; M = 11000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 183; Q = 60 }
    9  = @{ D = 44179; J = 15; K = 7000;  L = 7000;  M = 7000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 117; Q = 60 }
    10 = @{ D = 44186; J = 15; K = 7000;  L = 7000;  M = 7000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 117; Q = 60 }
    11 = @{ D = 44405; J = 45; K = 9000;  L = 9000;  M = 9000;  N = '$/caja 50 unidades'; O = 'Provincia de Quillota';        P = 180; Q = 50 }
    12 = @{ D = 45001; J = 40; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 167; Q = 60 }
    13 = @{ D = 44312; J = 30; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 167; Q = 60 }
    14 = @{ D = 45030; J = 50; K = 6000;  L = 6000;  M = 6000;  N = '$/caja 50 unidades'; O = 'Región de Arica y Parinacota'; P = 120; Q = 50 }
    15 = @{ D = 44315; J = 25; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';          P = 167; Q = 60 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q: Kg o Unidades
}
